# Update the calibration-service request list with the new batch of
# equipment records (replaces the previous sample rows with the new
# "NHP" equipment list) and refreshes the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Equipment"
$ws.Range("B1").Value = "Manufacturer"
$ws.Range("C1").Value = "Model"
$ws.Range("D1").Value = "Serial No."
$ws.Range("E1").Value = "ID No."
$ws.Range("F1").Value = "Calibration Date"
$ws.Range("G1").Value = "Certificate No."

# --- Data rows (row 2 - row 18) ---
# Row 2
$ws.Range("A2").Value = "HOT AIR OVEN"
$ws.Range("B2").Value = "BINDER"
$ws.Range("C2").Value = "RF 53"
$ws.Range("D2").Value = "RL15-07514"
$ws.Range("E2").Value = "NHP-CT-013"
$ws.Range("F2").Value = 45825
$ws.Range("G2").Value = "CH25052997"

# Row 3
$ws.Range("A3").Value = "CENTRIFUGE"
$ws.Range("B3").Value = "HETTICH"
$ws.Range("C3").Value = "ROTINA 46"
$ws.Range("D3").Value = "0000240"
$ws.Range("E3").Value = "NHP-C-CT-003"
$ws.Range("F3").Value = 45824
$ws.Range("G3").Value = "CF25052998"

# Row 4
$ws.Range("A4").Value = "ELECTRONIC BALANCE"
$ws.Range("B4").Value = "OHAUS"
$ws.Range("C4").Value = "V22PWE3T"
$ws.Range("D4").Value = "8338347025"
$ws.Range("E4").Value = "NHP-HIS-006"
$ws.Range("F4").Value = 45825
$ws.Range("G4").Value = "EB25052999"

# Row 5
$ws.Range("A5").Value = "ELECTRONIC BALANCE"
$ws.Range("B5").Value = "OHAUS"
$ws.Range("C5").Value = "FB3"
$ws.Range("D5").Value = "8027090157"
$ws.Range("E5").Value = "NHP-HIS-007"
$ws.Range("F5").Value = 45824
$ws.Range("G5").Value = "EB25053000"

# Row 6
$ws.Range("A6").Value = "ELECTRONIC BALANCE"
$ws.Range("B6").Value = "OHAUS"
$ws.Range("C6").Value = "V22PWE3T"
$ws.Range("D6").Value = "8337110661"
$ws.Range("E6").Value = "NHP-HIS-008"
$ws.Range("F6").Value = 45825
$ws.Range("G6").Value = "EB25053001"

# Row 7
$ws.Range("A7").Value = "ELECTRONIC BALANCE"
$ws.Range("B7").Value = "SARTORIUS"
$ws.Range("C7").Value = "BSA3202S-CW"
$ws.Range("D7").Value = "26290320"
$ws.Range("E7").Value = "NHP-HIS-009"
$ws.Range("F7").Value = 45824
$ws.Range("G7").Value = "EB25053002"

# Row 8
$ws.Range("A8").Value = "HOT AIR OVEN"
$ws.Range("B8").Value = "MEMMERT"
$ws.Range("C8").Value = "UF30"
$ws.Range("D8").Value = "B124.2173"
$ws.Range("E8").Value = "NHP-HIS-015"
$ws.Range("F8").Value = 45826
$ws.Range("G8").Value = "CH25053003"

# Row 9
$ws.Range("A9").Value = "HOT AIR OVEN"
$ws.Range("B9").Value = "MEMMERT"
$ws.Range("C9").Value = "UN30"
$ws.Range("D9").Value = "B121.0649"
$ws.Range("E9").Value = "NHP-HIS-016"
$ws.Range("F9").Value = 45825
$ws.Range("G9").Value = "CH25053004"

# Row 10
$ws.Range("A10").Value = "WATER BATH"
$ws.Range("B10").Value = "ZEEDO"
$ws.Range("C10").Value = "HS1125"
$ws.Range("D10").Value = "CJHS2099-32HS1125"
$ws.Range("E10").Value = "NHP-HIS-033"
$ws.Range("F10").Value = 45826
$ws.Range("G10").Value = "WB25053005"

# Row 11
$ws.Range("A11").Value = "WATER BATH"
$ws.Range("B11").Value = "ELECTROTHERMAL"
$ws.Range("C11").Value = "MH8517"
$ws.Range("D11").Value = "M553340/07"
$ws.Range("E11").Value = "NHP-HIS-034"
$ws.Range("F11").Value = 45825
$ws.Range("G11").Value = "WB25053006"

# Row 12
$ws.Range("A12").Value = "WATER BATH"
$ws.Range("B12").Value = "ELECTROTHERMAL"
$ws.Range("C12").Value = "MH8517"
$ws.Range("D12").Value = "M577610/08"
$ws.Range("E12").Value = "NHP-HIS-035"
$ws.Range("F12").Value = 45826
$ws.Range("G12").Value = "WB25053007"

# Row 13
$ws.Range("A13").Value = "WATER BATH"
$ws.Range("B13").Value = "LEICA"
$ws.Range("C13").Value = "HI1210"
$ws.Range("D13").Value = "16634"
$ws.Range("E13").Value = "NHP-HIS-036"
$ws.Range("F13").Value = 45825
$ws.Range("G13").Value = "WB25053008"

# Row 14
$ws.Range("A14").Value = "REFRIGERATOR"
$ws.Range("B14").Value = "ยี่ห้อ13"
$ws.Range("C14").Value = "รุ่น13"
$ws.Range("D14").Value = "หมายเลขเครื่อง13"
$ws.Range("E14").Value = "NHP-HIS-038"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "ยกเลิก(เครื่องเสีย)"

# Row 15
$ws.Range("A15").Value = "REFRIGERATOR"
$ws.Range("B15").Value = "MEDICOOL"
$ws.Range("C15").Value = "MDC145"
$ws.Range("D15").Value = "MDC145-202311046"
$ws.Range("E15").Value = "NHP-MOL-003"
$ws.Range("F15").Value = 45824
$ws.Range("G15").Value = "CH25053010"

# Row 16
$ws.Range("A16").Value = "FREEZER"
$ws.Range("B16").Value = "MEDICOOL"
$ws.Range("C16").Value = "MDC115"
$ws.Range("D16").Value = "MDC115-202311019"
$ws.Range("E16").Value = "NHP-MOL-004"
$ws.Range("F16").Value = 45824
$ws.Range("G16").Value = "CH25053011"

# Row 17
$ws.Range("A17").Value = "CENTRIFUGE"
$ws.Range("B17").Value = "LABTRON"
$ws.Range("C17").Value = "CF8"
$ws.Range("D17").Value = "1110801241A006"
$ws.Range("E17").Value = "NHP-MOL-006"
$ws.Range("F17").Value = 45826
$ws.Range("G17").Value = "CF25053012"

# Row 18
$ws.Range("A18").Value = "CENTRIFUGE"
$ws.Range("B18").Value = "UGAIYA"
$ws.Range("C18").Value = "TXD3"
$ws.Range("D18").Value = "W19800024030108"
$ws.Range("E18").Value = "NHP-CT-0048"
$ws.Range("F18").Value = 45824
$ws.Range("G18").Value = "CF25053009"

# --- Sheet view: zoom out to 70% (page break preview) and move the
#     active selection to G14 ---
$excel.ActiveWindow.View = 2
$ws.Range("G14").Select()
$excel.ActiveWindow.Zoom = 70
